# Updated cryptos list on Mon Jul 22 05:09:37 UTC 2024 with GitHub Actions
#
# Refreshes the price (column D) and 1h-volume-change (column E) figures
# for every coin row, and re-sorts a few rows whose relative ranking
# changed (Filecoin now above Maker, EnergySwap above Hedera, Bittensor
# above VeChain).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT, even when it looks like a number
# (e.g. "4.70", "0.999") so Excel's auto-type-detection doesn't silently
# turn it into a numeric cell (which would also drop trailing zeros).
function Set-TextValue {
    param($addr, $val)

    $cell = $ws.Range($addr)
    $looksNumeric = $val -match '^[+-]?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

# Row 2 - Bitcoin
Set-TextValue 'D2' '67.836.70'
Set-TextValue 'E2' '  +1.30%  '

# Row 3 - Ethereum
Set-TextValue 'D3' '3.510.22'
Set-TextValue 'E3' '  +0.29%  '

# Row 4 - TetherUSD
Set-TextValue 'E4' '  -0.01%  '

# Row 5 - BNB
Set-TextValue 'D5' '600.19'
Set-TextValue 'E5' '  +1.17%  '

# Row 6 - Solana
Set-TextValue 'D6' '180.79'
Set-TextValue 'E6' '  +4.62%  '

# Row 7 - USDC
Set-TextValue 'E7' '  +0.04%  '

# Row 8 - LidoStakedEther
Set-TextValue 'D8' '3.509.32'
Set-TextValue 'E8' '  +0.28%  '

# Row 9 - XRP
Set-TextValue 'D9' '0.594'
Set-TextValue 'E9' '  -1.65%  '

# Row 10 - Dogecoin
Set-TextValue 'E10' '  +7.63%  '

# Row 11 - Toncoin
Set-TextValue 'D11' '7.15'
Set-TextValue 'E11' '  -1.65%  '

# Row 12 - Cardano
Set-TextValue 'D12' '0.436'
Set-TextValue 'E12' '  +0.83%  '

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue 'D13' '4.120.98'
Set-TextValue 'E13' '  +0.33%  '

# Row 14 - Avalanche
Set-TextValue 'D14' '32.60'
Set-TextValue 'E14' '  +12.83%  '

# Row 15 - TRON
Set-TextValue 'E15' '  +1.34%  '

# Row 16 - WrappedBTC
Set-TextValue 'D16' '67.855.88'
Set-TextValue 'E16' '  +1.33%  '

# Row 17 - ShibaInu
Set-TextValue 'E17' '  +0.65%  '

# Row 18 - WrappedEther
Set-TextValue 'D18' '3.511.53'
Set-TextValue 'E18' '  +0.94%  '

# Row 19 - Polkadot
Set-TextValue 'D19' '6.33'
Set-TextValue 'E19' '  +0.42%  '

# Row 20 - Chainlink
Set-TextValue 'D20' '14.48'
Set-TextValue 'E20' '  +2.40%  '

# Row 21 - BitcoinCash
Set-TextValue 'D21' '398.85'
Set-TextValue 'E21' '  +1.48%  '

# Row 22 - Uniswap
Set-TextValue 'D22' '7.96'
Set-TextValue 'E22' '  -0.22%  '

# Row 23 - Litecoin
Set-TextValue 'D23' '73.62'
Set-TextValue 'E23' '  +0.77%  '

# Row 24 - Polygon
Set-TextValue 'D24' '0.541'
Set-TextValue 'E24' '  +1.07%  '

# Row 25 - Dai
Set-TextValue 'D25' '0.999'
Set-TextValue 'E25' '  -0.03%  '

# Row 26 - LEO
Set-TextValue 'E26' '  +0.57%  '

# Row 27 - PEPE
Set-TextValue 'E27' '  +1.83%  '

# Row 28 - InternetComputer(DFINITY)
Set-TextValue 'D28' '10.55'
Set-TextValue 'E28' '  +3.42%  '

# Row 29 - Kaspa
Set-TextValue 'E29' '  -2.31%  '

# Row 30 - Binance-PegBSC-USD
Set-TextValue 'D30' '0.998'
Set-TextValue 'E30' '  +0.05%  '

# Row 31 - NEARProtocol
Set-TextValue 'D31' '6.25'
Set-TextValue 'E31' '  -0.22%  '

# Row 32 - Fetch.AI
Set-TextValue 'D32' '1.45'
Set-TextValue 'E32' '  +0.52%  '

# Row 33 - PancakeSwap
Set-TextValue 'E33' '  +1.68%  '

# Row 34 - EthereumClassic
Set-TextValue 'D34' '23.85'
Set-TextValue 'E34' '  +0.50%  '

# Row 35 - Aptos
Set-TextValue 'D35' '7.50'
Set-TextValue 'E35' '  +2.22%  '

# Row 37 - ImmutableX
Set-TextValue 'E37' '  -2.59%  '

# Row 38 - Monero
Set-TextValue 'D38' '164.84'
Set-TextValue 'E38' '  +0.87%  '

# Row 39 - Mantle
Set-TextValue 'D39' '0.878'
Set-TextValue 'E39' '  -0.29%  '

# Row 40 - Stacks
Set-TextValue 'E40' '  +1.35%  '

# Row 41 - dogwifhat
Set-TextValue 'D41' '2.77'
Set-TextValue 'E41' '  +9.10%  '

# Row 42 - RenderToken
Set-TextValue 'D42' '6.97'
Set-TextValue 'E42' '  -0.29%  '

# Row 43 - was Maker, now Filecoin (rows 43/44 swap ranking)
Set-TextValue 'B43' 'Filecoin'
Set-TextValue 'C43' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D43' '4.70'
Set-TextValue 'E43' '  +0.71%  '

# Row 44 - was Filecoin, now Maker
Set-TextValue 'B44' 'Maker'
Set-TextValue 'C44' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D44' '2.876.25'
Set-TextValue 'E44' '  +2.45%  '

# Row 45 - was Hedera, now EnergySwap (rows 45/46 swap ranking)
Set-TextValue 'B45' 'EnergySwap'
Set-TextValue 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '26.26'
Set-TextValue 'E45' '  -0.30%  '

# Row 46 - was EnergySwap, now Hedera
Set-TextValue 'B46' 'Hedera'
Set-TextValue 'C46' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D46' '0.0733'
Set-TextValue 'E46' '  -1.61%  '

# Row 47 - InjectiveProtocol
Set-TextValue 'D47' '26.73'
Set-TextValue 'E47' '  -2.30%  '

# Row 48 - OKB
Set-TextValue 'D48' '42.29'
Set-TextValue 'E48' '  -0.96%  '

# Row 49 - was VeChain, now Bittensor (rows 49/50 swap ranking)
Set-TextValue 'B49' 'Bittensor'
Set-TextValue 'C49' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D49' '345.75'
Set-TextValue 'E49' '  +2.70%  '

# Row 50 - was Bittensor, now VeChain
Set-TextValue 'B50' 'VeChain'
Set-TextValue 'C50' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D50' '0.0302'
Set-TextValue 'E50' '  +0.25%  '

# Row 51 - ONDO
Set-TextValue 'E51' '  -0.85%  '
